$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "261.73";  "E2"  = "0.89%"
    "D3"  = "27.12";   "E3"  = "0.40%"
    "D4"  = "4.699"
    "E5"  = "2.84%"
    "D6"  = "6.730";   "E6"  = "0.93%"
    "D7"  = "0.8500";  "E7"  = "-1.05%"
    "D8"  = "0.9195";  "E8"  = "-1.35%"
    "E9"  = "1.15%"
    "D10" = "0.04541"; "E10" = "-1.58%"
    "D11" = "0.07089"; "E11" = "1.18%"
    "D12" = "0.03131"; "E12" = "0.47%"
    "D13" = "0.09045"; "E13" = "-0.99%"
    "D14" = "0.001530";"E14" = "-0.10%"
    "D15" = "0.0006150";"E15" = "1.35%"
    "D16" = "0.006026";"E16" = "-1.74%"
    "D17" = "3.466";   "E17" = "0.02%"
    "D18" = "3.164";   "E18" = "0.20%"
    "D19" = "2.163";   "E19" = "-0.12%"
    "E20" = "-0.23%"
    "E21" = "1.00%"
    "D22" = "4.099";   "E22" = "-0.84%"
    "D23" = "0.04243"; "E23" = "0.23%"
    "E24" = "0.27%"
    "E25" = "-6.12%"
    "E26" = "0.29%"
    "D40" = "0.03929"; "E40" = "2.11%"
    "D41" = "0.1114";  "E41" = "-0.29%"
    "D42" = "0.004133";"E42" = "6.04%"
    "E43" = "-9.54%"
    "D44" = "0.01383"; "E44" = "-9.43%"
    "E45" = "1.20%"
    "E46" = "0.27%"
    "E47" = "-28.15%"
    "D48" = "0.1677";  "E48" = "28.72%"
    "D49" = "0.00002101";"E49" = "0.27%"
    "D50" = "0.0002001";"E50" = "0.27%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
